$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2025
$ws.Range("J18").Value = 2050
$ws.Range("L18").Value = 2050
$ws.Range("N18").Value = -2618

$ws.Range("H92").Value = 21739446
$ws.Range("I92").Value = 22727538
$ws.Range("K92").Value = 22727538
$ws.Range("M92").Value = -22726290

$ws.Range("H99").Value = 852.3077
$ws.Range("I99").Value = 565.5
$ws.Range("J99").Value = 1311.2
$ws.Range("K99").Value = 1696.5
$ws.Range("L99").Value = 3933.6
$ws.Range("M99").Value = -198.5
$ws.Range("N99").Value = -6929.6

$ws.Range("H135").Value = 1193.5
$ws.Range("I135").Value = 881.6667
$ws.Range("K135").Value = 7935.0003
$ws.Range("M135").Value = -5400.0003

$ws.Range("H138").Value = 3218.9473
$ws.Range("I138").Value = 1970
$ws.Range("J138").Value = 3727.7778
$ws.Range("K138").Value = 5910
$ws.Range("L138").Value = 11183.3334
$ws.Range("M138").Value = -770
$ws.Range("N138").Value = -21463.3334

$ws.Range("H141").Value = 4820.385
$ws.Range("I141").Value = 4697.273
$ws.Range("K141").Value = 14091.819
$ws.Range("M141").Value = -8911.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 503
$ws.Range("I14").Value = 6
$ws.Range("K14").Value = 6
$ws.Range("M14").Value = 169

$ws.Range("H32").Value = 3501.2307
$ws.Range("I32").Value = 2686.6858
$ws.Range("K32").Value = 2686.6858
$ws.Range("M32").Value = -2399.6858

$ws.Range("H43").Value = 40000
$ws.Range("J43").Value = 40000
$ws.Range("L43").Value = 40000
$ws.Range("N43").Value = -40626

$ws.Range("H122").Value = 5250.1924
$ws.Range("I122").Value = 3531.9375
$ws.Range("J122").Value = 7999.4
$ws.Range("K122").Value = 10595.8125
$ws.Range("L122").Value = 23998.2
$ws.Range("M122").Value = -8145.8125
$ws.Range("N122").Value = -28898.2

$ws.Range("H132").Value = 3076.4517
$ws.Range("I132").Value = 2845.1482
$ws.Range("K132").Value = 8535.444600000001
$ws.Range("M132").Value = -6005.444600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 667.5714
$ws.Range("I86").Value = 498
$ws.Range("J86").Value = 695.8333
$ws.Range("K86").Value = 498
$ws.Range("L86").Value = 695.8333
$ws.Range("M86").Value = 625
$ws.Range("N86").Value = -2941.8333

$ws.Range("H89").Value = 667.5714
$ws.Range("I89").Value = 498
$ws.Range("J89").Value = 695.8333
$ws.Range("K89").Value = 2490
$ws.Range("L89").Value = 3479.1665
$ws.Range("M89").Value = 3126
$ws.Range("N89").Value = -14711.1665

$ws.Range("H134").Value = 10754699
$ws.Range("I134").Value = 1685.3846
$ws.Range("K134").Value = 5056.1538
$ws.Range("M134").Value = -2521.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 858.8
$ws.Range("I22").Value = 898.5
$ws.Range("J22").Value = 700
$ws.Range("K22").Value = 898.5
$ws.Range("L22").Value = 700
$ws.Range("M22").Value = -548.5
$ws.Range("N22").Value = -1400

$ws.Range("H132").Value = 5150.875
$ws.Range("J132").Value = 6000
$ws.Range("L132").Value = 18000
$ws.Range("N132").Value = -23060

$ws.Range("H141").Value = 392018.7
$ws.Range("J141").Value = 392018.7
$ws.Range("L141").Value = 392018.7
$ws.Range("N141").Value = -402378.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 330.42856
$ws.Range("I103").Value = 350
$ws.Range("K103").Value = 1050
$ws.Range("M103").Value = -171

$ws.Range("H131").Value = 1652.2954
$ws.Range("J131").Value = 1920.0344
$ws.Range("L131").Value = 5760.1032
$ws.Range("N131").Value = -15840.1032

$ws.Range("H132").Value = 603.25
$ws.Range("I132").Value = 603.25
$ws.Range("K132").Value = 5429.25
$ws.Range("M132").Value = -2899.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 46912.5
$ws.Range("I62").Value = 44000
$ws.Range("J62").Value = 49825
$ws.Range("K62").Value = 44000
$ws.Range("L62").Value = 49825
$ws.Range("M62").Value = -43314
$ws.Range("N62").Value = -51197

$ws.Range("H65").Value = 46912.5
$ws.Range("I65").Value = 44000
$ws.Range("J65").Value = 49825
$ws.Range("K65").Value = 132000
$ws.Range("L65").Value = 149475
$ws.Range("M65").Value = -128568
$ws.Range("N65").Value = -156339

$ws.Range("H80").Value = 3115.5
$ws.Range("I80").Value = 2898
$ws.Range("J80").Value = 3224.25
$ws.Range("K80").Value = 2898
$ws.Range("L80").Value = 3224.25
$ws.Range("M80").Value = -1900
$ws.Range("N80").Value = -5220.25

$ws.Range("H83").Value = 3115.5
$ws.Range("I83").Value = 2898
$ws.Range("J83").Value = 3224.25
$ws.Range("K83").Value = 14490
$ws.Range("L83").Value = 16121.25
$ws.Range("M83").Value = -9498
$ws.Range("N83").Value = -26105.25

$ws.Range("H122").Value = 2013.7858
$ws.Range("I122").Value = 1698.8334
$ws.Range("K122").Value = 5096.5002
$ws.Range("M122").Value = -2646.5002

$ws.Range("H132").Value = 2534.7083
$ws.Range("I132").Value = 2532.8635
$ws.Range("K132").Value = 7598.5905
$ws.Range("M132").Value = -5068.5905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5500.857
$ws.Range("I22").Value = 4624
$ws.Range("J22").Value = 5851.6
$ws.Range("K22").Value = 4624
$ws.Range("L22").Value = 5851.6
$ws.Range("M22").Value = -4329
$ws.Range("N22").Value = -6441.6

$ws.Range("H27").Value = 5500.857
$ws.Range("I27").Value = 4624
$ws.Range("J27").Value = 5851.6
$ws.Range("K27").Value = 4624
$ws.Range("L27").Value = 5851.6
$ws.Range("M27").Value = -4517
$ws.Range("N27").Value = -6065.6

$ws.Range("H40").Value = 4088.9092
$ws.Range("I40").Value = 3418.7896
$ws.Range("K40").Value = 3418.7896
$ws.Range("M40").Value = -3282.7896

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = 0

$ws.Range("H68").Value = 2869.5715
$ws.Range("I68").Value = 2653.0908
$ws.Range("K68").Value = 2653.0908
$ws.Range("M68").Value = -1904.0908

$ws.Range("H71").Value = 2869.5715
$ws.Range("I71").Value = 2653.0908
$ws.Range("K71").Value = 13265.454
$ws.Range("M71").Value = -9521.454

$ws.Range("H82").Value = 1846.0526
$ws.Range("I82").Value = 1823.375
$ws.Range("J82").Value = 1862.5454
$ws.Range("K82").Value = 1823.375
$ws.Range("L82").Value = 1862.5454
$ws.Range("M82").Value = -1462.375
$ws.Range("N82").Value = -2584.5454

$ws.Range("H85").Value = 1846.0526
$ws.Range("I85").Value = 1823.375
$ws.Range("J85").Value = 1862.5454
$ws.Range("K85").Value = 1823.375
$ws.Range("L85").Value = 1862.5454
$ws.Range("M85").Value = -575.375
$ws.Range("N85").Value = -4358.5454

$ws.Range("H93").Value = 111112440
$ws.Range("J93").Value = 1972
$ws.Range("L93").Value = 1972
$ws.Range("N93").Value = -4468

$ws.Range("H112").Value = 127499
$ws.Range("J112").Value = 127499
$ws.Range("L112").Value = 127499
$ws.Range("N112").Value = -130453

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 26678
$ws.Range("I21").Value = 20000
$ws.Range("K21").Value = 20000
$ws.Range("M21").Value = -19765

$ws.Range("H35").Value = 26678
$ws.Range("I35").Value = 20000
$ws.Range("K35").Value = 20000
$ws.Range("M35").Value = -19710

$ws.Range("H122").Value = 5952
$ws.Range("J122").Value = 13230.2
$ws.Range("L122").Value = 39690.60000000001
$ws.Range("N122").Value = -44590.60000000001

$ws.Range("H132").Value = 3476.3704
$ws.Range("I132").Value = 3217.3489
$ws.Range("K132").Value = 9652.046699999999
$ws.Range("M132").Value = -7122.046699999999

$ws.Range("H136").Value = 1755.228
$ws.Range("I136").Value = 1242.3721
$ws.Range("J136").Value = 3330.4285
$ws.Range("K136").Value = 3727.1163
$ws.Range("L136").Value = 9991.2855
$ws.Range("M136").Value = -1177.1163
$ws.Range("N136").Value = -15091.2855
